# Registro de pago temporal parcial
# Update the "Casos de Uso" sheet: task F11 status moves from "Por iniciar"
# to "En proceso", and 2 hours of consumption are logged against it in
# column Q (day 4), which ripples the running "Rest."/"Cons." formulas
# already present in the sheet. Also restore the active-cell selection in
# the frozen bottom-right pane to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Task status: "Por iniciar" -> "En proceso"
$ws.Range("F11").Value = "En proceso"

# Log 2 hours consumed in the day-4 column; downstream Rest./Total formulas
# recompute automatically.
$ws.Range("Q11").Value = 2

# Move the active selection (bottom-right frozen pane) to F12.
$ws.Range("F12").Select()
